# Apply the edit described by the diff:
# 1) Update the '总计' summary sheet with a new 2022-Q3 row (shifting the rest down).
# 2) Insert a new '2022-Q3' worksheet (right after '总计') holding the fund-holdings table.

$wb = $excel.ActiveWorkbook

$summary = $wb.Worksheets.Item("总计")

# --- Step 1: grow the summary table by one row, preserving row formatting ---
# Duplicate the last data row (row 7) into the new row 8 so the new row inherits
# the same cell styling (bold/bordered/centered index column, etc.)
$summary.Range("A7:D7").Copy($summary.Range("A8:D8"))

# --- Step 2: (re)write every data row A2:D8 with the post-edit values ---
$summary.Cells.Item(2, 1).Value = 0
$summary.Cells.Item(2, 2).Value = '2022-Q3'
$summary.Cells.Item(2, 3).Value = '13'
$summary.Cells.Item(2, 4).Value = '3.11'

$summary.Cells.Item(3, 1).Value = 1
$summary.Cells.Item(3, 2).Value = '2022-Q1'
$summary.Cells.Item(3, 3).Value = '22'
$summary.Cells.Item(3, 4).Value = '8.67'

$summary.Cells.Item(4, 1).Value = 2
$summary.Cells.Item(4, 2).Value = '2021-Q4'
$summary.Cells.Item(4, 3).Value = '17'
$summary.Cells.Item(4, 4).Value = '6.51'

$summary.Cells.Item(5, 1).Value = 3
$summary.Cells.Item(5, 2).Value = '2021-Q3'
$summary.Cells.Item(5, 3).Value = '24'
$summary.Cells.Item(5, 4).Value = '12.35'

$summary.Cells.Item(6, 1).Value = 4
$summary.Cells.Item(6, 2).Value = '2021-Q2'
$summary.Cells.Item(6, 3).Value = '23'
$summary.Cells.Item(6, 4).Value = '12.08'

$summary.Cells.Item(7, 1).Value = 5
$summary.Cells.Item(7, 2).Value = '2021-Q1'
$summary.Cells.Item(7, 3).Value = '19'
$summary.Cells.Item(7, 4).Value = '14.67'

$summary.Cells.Item(8, 1).Value = 6
$summary.Cells.Item(8, 2).Value = '2020-Q4'
$summary.Cells.Item(8, 3).Value = '17'
$summary.Cells.Item(8, 4).Value = '4.29'

# --- Step 3: insert the new '2022-Q3' worksheet right after '总计' ---
$newWs = $wb.Worksheets.Add($null, $summary)
$newWs.Name = "2022-Q3"

# Header row (B1:H1)
$newWs.Cells.Item(1, 2).Value = '基金代码'
$newWs.Cells.Item(1, 3).Value = '基金名称'
$newWs.Cells.Item(1, 4).Value = '基金规模'
$newWs.Cells.Item(1, 5).Value = '股票总仓位'
$newWs.Cells.Item(1, 6).Value = '仓位占比'
$newWs.Cells.Item(1, 7).Value = '持有市值(亿元)'
$newWs.Cells.Item(1, 8).Value = '仓位排名'

# Data rows A2:H14
# Row 2
$newWs.Cells.Item(2, 1).Value = 0
$newWs.Cells.Item(2, 2).Value = "'" + '513090'
$newWs.Cells.Item(2, 3).Value = '易方达中证香港证券投资主题ETF'
$newWs.Cells.Item(2, 4).Value = "'" + '10.53'
$newWs.Cells.Item(2, 5).Value = "'" + '96.33'
$newWs.Cells.Item(2, 6).Value = "'" + '10.68'
$newWs.Cells.Item(2, 7).Value = "'" + '1.1246'
$newWs.Cells.Item(2, 8).Value = 3

# Row 3
$newWs.Cells.Item(3, 1).Value = 1
$newWs.Cells.Item(3, 2).Value = "'" + '014362'
$newWs.Cells.Item(3, 3).Value = '睿远稳进配置两年持有混合A'
$newWs.Cells.Item(3, 4).Value = "'" + '64.40'
$newWs.Cells.Item(3, 5).Value = "'" + '35.09'
$newWs.Cells.Item(3, 6).Value = "'" + '0.97'
$newWs.Cells.Item(3, 7).Value = "'" + '0.6247'
$newWs.Cells.Item(3, 8).Value = 10

# Row 4
$newWs.Cells.Item(4, 1).Value = 2
$newWs.Cells.Item(4, 2).Value = "'" + '014363'
$newWs.Cells.Item(4, 3).Value = '睿远稳进配置两年持有混合C'
$newWs.Cells.Item(4, 4).Value = "'" + '35.67'
$newWs.Cells.Item(4, 5).Value = "'" + '35.09'
$newWs.Cells.Item(4, 6).Value = "'" + '0.97'
$newWs.Cells.Item(4, 7).Value = "'" + '0.3460'
$newWs.Cells.Item(4, 8).Value = 10

# Row 5
$newWs.Cells.Item(5, 1).Value = 3
$newWs.Cells.Item(5, 2).Value = "'" + '012943'
$newWs.Cells.Item(5, 3).Value = '广发稳睿六个月持有期混合A'
$newWs.Cells.Item(5, 4).Value = "'" + '20.75'
$newWs.Cells.Item(5, 5).Value = "'" + '26.11'
$newWs.Cells.Item(5, 6).Value = "'" + '1.50'
$newWs.Cells.Item(5, 7).Value = "'" + '0.3112'
$newWs.Cells.Item(5, 8).Value = 6

# Row 6
$newWs.Cells.Item(6, 1).Value = 4
$newWs.Cells.Item(6, 2).Value = "'" + '012944'
$newWs.Cells.Item(6, 3).Value = '广发稳睿六个月持有期混合C'
$newWs.Cells.Item(6, 4).Value = "'" + '19.14'
$newWs.Cells.Item(6, 5).Value = "'" + '26.11'
$newWs.Cells.Item(6, 6).Value = "'" + '1.50'
$newWs.Cells.Item(6, 7).Value = "'" + '0.2871'
$newWs.Cells.Item(6, 8).Value = 6

# Row 7
$newWs.Cells.Item(7, 1).Value = 5
$newWs.Cells.Item(7, 2).Value = "'" + '202801'
$newWs.Cells.Item(7, 3).Value = '南方全球精选配置（QDII-FOF）'
$newWs.Cells.Item(7, 4).Value = "'" + '15.80'
$newWs.Cells.Item(7, 5).Value = "'" + '29.52'
$newWs.Cells.Item(7, 6).Value = "'" + '1.31'
$newWs.Cells.Item(7, 7).Value = "'" + '0.2070'
$newWs.Cells.Item(7, 8).Value = 7

# Row 8
$newWs.Cells.Item(8, 1).Value = 6
$newWs.Cells.Item(8, 2).Value = "'" + '007109'
$newWs.Cells.Item(8, 3).Value = '南方沪港深核心优势混合'
$newWs.Cells.Item(8, 4).Value = "'" + '1.59'
$newWs.Cells.Item(8, 5).Value = "'" + '85.13'
$newWs.Cells.Item(8, 6).Value = "'" + '5.20'
$newWs.Cells.Item(8, 7).Value = "'" + '0.0827'
$newWs.Cells.Item(8, 8).Value = 4

# Row 9
$newWs.Cells.Item(9, 1).Value = 7
$newWs.Cells.Item(9, 2).Value = "'" + '011355'
$newWs.Cells.Item(9, 3).Value = '华泰柏瑞港股通时代机遇混合A'
$newWs.Cells.Item(9, 4).Value = "'" + '0.54'
$newWs.Cells.Item(9, 5).Value = "'" + '91.80'
$newWs.Cells.Item(9, 6).Value = "'" + '6.67'
$newWs.Cells.Item(9, 7).Value = "'" + '0.0360'
$newWs.Cells.Item(9, 8).Value = 7

# Row 10
$newWs.Cells.Item(10, 1).Value = 8
$newWs.Cells.Item(10, 2).Value = "'" + '005576'
$newWs.Cells.Item(10, 3).Value = '华泰柏瑞新金融地产灵活配置混合A'
$newWs.Cells.Item(10, 4).Value = "'" + '0.84'
$newWs.Cells.Item(10, 5).Value = "'" + '94.39'
$newWs.Cells.Item(10, 6).Value = "'" + '3.44'
$newWs.Cells.Item(10, 7).Value = "'" + '0.0289'
$newWs.Cells.Item(10, 8).Value = 8

# Row 11
$newWs.Cells.Item(11, 1).Value = 9
$newWs.Cells.Item(11, 2).Value = "'" + '003413'
$newWs.Cells.Item(11, 3).Value = '华泰柏瑞新经济沪港深混合'
$newWs.Cells.Item(11, 4).Value = "'" + '0.42'
$newWs.Cells.Item(11, 5).Value = "'" + '86.45'
$newWs.Cells.Item(11, 6).Value = "'" + '6.21'
$newWs.Cells.Item(11, 7).Value = "'" + '0.0261'
$newWs.Cells.Item(11, 8).Value = 2

# Row 12
$newWs.Cells.Item(12, 1).Value = 10
$newWs.Cells.Item(12, 2).Value = "'" + '460010'
$newWs.Cells.Item(12, 3).Value = '华泰柏瑞亚洲领导企业混合（QDII）'
$newWs.Cells.Item(12, 4).Value = "'" + '0.36'
$newWs.Cells.Item(12, 5).Value = "'" + '93.44'
$newWs.Cells.Item(12, 6).Value = "'" + '5.70'
$newWs.Cells.Item(12, 7).Value = "'" + '0.0205'
$newWs.Cells.Item(12, 8).Value = 7

# Row 13
$newWs.Cells.Item(13, 1).Value = 11
$newWs.Cells.Item(13, 2).Value = "'" + '011356'
$newWs.Cells.Item(13, 3).Value = '华泰柏瑞港股通时代机遇混合C'
$newWs.Cells.Item(13, 4).Value = "'" + '0.24'
$newWs.Cells.Item(13, 5).Value = "'" + '91.80'
$newWs.Cells.Item(13, 6).Value = "'" + '6.67'
$newWs.Cells.Item(13, 7).Value = "'" + '0.0160'
$newWs.Cells.Item(13, 8).Value = 7

# Row 14
$newWs.Cells.Item(14, 1).Value = 12
$newWs.Cells.Item(14, 2).Value = "'" + '016374'
$newWs.Cells.Item(14, 3).Value = '华泰柏瑞新金融地产灵活配置混合C'
$newWs.Cells.Item(14, 4).Value = "'" + '0.02'
$newWs.Cells.Item(14, 5).Value = "'" + '94.39'
$newWs.Cells.Item(14, 6).Value = "'" + '3.44'
$newWs.Cells.Item(14, 7).Value = "'" + '0.0007'
$newWs.Cells.Item(14, 8).Value = 8

